$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("batiment")

# Insert a new row at position 10, shifting existing rows 10-17 down to 11-18
$ws.Rows.Item(10).Insert()

# Fill in the new row 10 with the code_epci entry
$ws.Range("A10").Value = "code_epci"
$ws.Range("B10").Value = "TEXT"
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = "Identifiant EPCI (établissements publics de coopération intercommunale) INSEE"
